# trial_info.xlsx update — add two new video-frame reference columns
# (video_underwater_frame_catch, video_underwater_frame_all_cups_in_water)
# between the existing "farfield split seconds" data and the
# video_above_water columns, fill in new data for throws on 4/27, and add
# a small scratch area (rows 12-15) computing frame-time offsets (crude
# look at where the beamwidth change came from).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Shift the old L/M ("video_above_water_filename" /
#    "video_above_water_frame_throw") columns right by two, to N/O, and
#    populate the new L/M columns with the catch / all-cups-in-water
#    frame data.
# ---------------------------------------------------------------------

# NOTE: the COM `.Value` *getter* does not return real cell data in this
# host (it echoes a property descriptor string) - `.Value2` reads back
# correctly, so every read below uses `.Value2` while every write still
# uses `.Value` (the setter works fine).

# -- Header row --------------------------------------------------------
$ws.Range("N1").Value = $ws.Range("L1").Value2
$ws.Range("O1").Value = $ws.Range("M1").Value2
$ws.Range("L1").Value = "video_underwater_frame_all_cups_in_water"
$ws.Range("M1").Value = "video_underwater_frame_catch"

# -- Row 5 (trial 1, throw 1, pm session) -------------------------------
$ws.Range("N5").Value = $ws.Range("L5").Value2
$ws.Range("O5").Value = $ws.Range("M5").Value2
$ws.Range("L5").ClearContents()
$ws.Range("M5").ClearContents()

# -- Row 6 ---------------------------------------------------------------
$ws.Range("N6").Value = $ws.Range("L6").Value2
$ws.Range("O6").Value = $ws.Range("M6").Value2
$ws.Range("L6").ClearContents()
$ws.Range("M6").ClearContents()

# -- Row 7 -----------------------------------------------------------
$ws.Range("N7").Value = $ws.Range("L7").Value2
$ws.Range("O7").Value = $ws.Range("M7").Value2
$ws.Range("K7").Value = 7883
$ws.Range("L7").Value = 7987
$ws.Range("M7").Value = 8097

# -- Row 8 -----------------------------------------------------------
$ws.Range("N8").Value = $ws.Range("L8").Value2
$ws.Range("O8").Value = $ws.Range("M8").Value2
$ws.Range("L8").Value = 9909
$ws.Range("M8").ClearContents()

# -- Row 9 -----------------------------------------------------------
$ws.Range("N9").Value = $ws.Range("L9").Value2
$ws.Range("O9").Value = $ws.Range("M9").Value2
$ws.Range("L9").Value = 11458
$ws.Range("M9").ClearContents()

# -- Row 10 -----------------------------------------------------------
$ws.Range("N10").Value = $ws.Range("L10").Value2
$ws.Range("O10").Value = $ws.Range("M10").Value2
$ws.Range("L10").Value = 14445
$ws.Range("M10").ClearContents()

# -- Row 2: new "video_underwater_frame_catch" value, no shift needed
#    because this row never had any L/M data.
$ws.Range("M2").Value = 9376

# ---------------------------------------------------------------------
# 2. Column widths: columns K..M are now one contiguous 35-wide block,
#    N keeps the 26-wide look, O keeps the ~29.8-wide look.
# ---------------------------------------------------------------------
$ws.Range("K1:M1").ColumnWidth = 34.166666666666664   # -> stored width 35
$ws.Range("N1").ColumnWidth = 25.166666666666668       # -> stored width 26
$ws.Range("O1").ColumnWidth = 28.998697916666668        # -> stored width ~29.83

# ---------------------------------------------------------------------
# 3. New scratch rows 12-15: frame-number -> seconds conversions
#    (crude look at where the beamwidth change came from).
# ---------------------------------------------------------------------
$ws.Range("J12").Formula = "=J7/29.97"
$ws.Range("K12").Formula = "=K7/29.97"
$ws.Range("M12").Formula = "=M7/29.97"

$ws.Range("K13").Formula = "=K12-J12"
$ws.Range("M13").Formula = "=M12-K12"

$ws.Range("J14").Formula = "=J2/29.97"
$ws.Range("K14").Formula = "=K2/29.97"
$ws.Range("M14").Formula = "=M2/29.97"

$ws.Range("K15").Formula = "=K14-J14"
$ws.Range("M15").Formula = "=M14-K14"

# ---------------------------------------------------------------------
# 4. View state: scroll right so column F is leftmost, select L11.
# ---------------------------------------------------------------------
$ws.Range("L11").Select()
